# edit.ps1 -- applies the "My Anime List Clone" project proposal revision
# Strategy: for each affected paragraph, replace the paragraph's run content
# (but not its paragraph mark, so pPr/numbering stays intact) with the exact
# target WordprocessingML via Range.InsertXML. This lets us place <w:proofErr>
# spell-check markers and relocate the <w:bookmarkStart>/<w:bookmarkEnd>
# "_GoBack" bookmark exactly where the target XML puts them.

function Set-ParaInnerXml {
    param(
        $doc,
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $para = $doc.Paragraphs($ParaIndex)
    $full = $para.Range
    # Shrink by 1 char so the trailing paragraph-mark (and thus <w:pPr>) is
    # left alone; only the run content inside the paragraph gets replaced.
    $target = $doc.Range($full.Start, $full.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $InnerXml +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

$d = $word.ActiveDocument

# --- 1. "...favorite animes, perhaps..." -> split out "animes" w/ proofErr ---
$p7 = '<w:r><w:t xml:space="preserve">The API should have an extensive amount of anime listed with good details on each, summary, genres, maybe ratings etc. Users can make lists of their favorite </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>animes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, perhaps add their own ratings to the lists. Users will also be able to share their lists, not necessarily with other users. </w:t></w:r>'

# --- 2. "...id, username, password, email, avatar" -> drop the _GoBack bookmark ---
$p10 = '<w:r><w:t xml:space="preserve">There should be a user with: id, username, password, email, </w:t></w:r><w:r><w:t>avatar</w:t></w:r>'

# --- 3. "Add_anime_to_list: id, list_id, user_id, anime api id, ..." -> split ids w/ proofErr ---
$p12 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Add_anime_to_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: id, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>list_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, anime </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> id, user’s personal rating, user’s comments on </w:t></w:r><w:r><w:t>anime choice</w:t></w:r>'

# --- 4. "API may be down ... I could possibly ..." -> move _GoBack bookmark mid-sentence ---
$p14 = '<w:r><w:tab/></w:r><w:r><w:t>API may be down in which case info on anime may not display, also wouldn’t be able to get recommendations for new anime, I co</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">uld possibly make a page to display if the API is down saying ‘having some technical issues’ or something to that affect as there won’t be much to show on my application if the API is down. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>'

# --- 5. "...favorite animes, ... more animes ... which animes fall..." -> split w/ proofErr ---
$p20 = '<w:r><w:t xml:space="preserve">User would need to register to make their own lists, but otherwise lists are public. Once registered they will have a portal from which to make list of favorite </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>animes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, they can add comments and ratings, remove and add more </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>animes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from each list, if not the user that owns the list they will not be able to make changes to the list. There will also be a page with a list of genres to see which </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>animes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fall into the genres they chose. Users can also click on the anime to be redirected to a page with more info on an anime. </w:t></w:r>'

# --- 6. "Eventually there should be full search functionality... " -> reorder/replace sentences ---
$p22 = '<w:r><w:t xml:space="preserve">Maybe an option for users to directly send recommendations for each other, or if they see a list is lacking, they can submit to a recommendation tab on the individual list. They should also be able to leave comments on lists or underneath a user’s comments. </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Maybe an alert sent to a user’s email when there is a new recommendation. Alerts on the app itself for commenting. </w:t></w:r><w:r><w:t xml:space="preserve">Will look more deeply into recommendation algorithms and see if there is a way to incorporate those into the application, as right now the recommendations only come from a set determined by the API, but not sure what it’s based on. </w:t></w:r>'

Set-ParaInnerXml $d 7  $p7
Set-ParaInnerXml $d 10 $p10
Set-ParaInnerXml $d 12 $p12
Set-ParaInnerXml $d 14 $p14
Set-ParaInnerXml $d 20 $p20
Set-ParaInnerXml $d 22 $p22
